$wb = $excel.ActiveWorkbook

# --- Update the conversion-of-the-day note on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.6 = 59868.76 pesos`n✅ 59868.76 pesos = 14.55 = 971.9 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update the rate cells on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 68.5
$ws2.Range("O10").Value = 4101.01
$ws2.Range("N12").Value = 4114
$ws2.Range("O12").Value = 66.786
